$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking strings
# (e.g. "197.04", "0.610") are preserved as text instead of being
# auto-converted to numbers by Excel.

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = '69.418.84'
$r = $ws.Range("E2")
$r.NumberFormat = "@"
$r.Value = '  -1.05%  '
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = '3.542.37'
$r = $ws.Range("E3")
$r.NumberFormat = "@"
$r.Value = '  -1.66%  '
$r = $ws.Range("E4")
$r.NumberFormat = "@"
$r.Value = '  +0.09%  '
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '197.04'
$r = $ws.Range("E5")
$r.NumberFormat = "@"
$r.Value = '  +0.23%  '
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '584.08'
$r = $ws.Range("E6")
$r.NumberFormat = "@"
$r.Value = '  -3.24%  '
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = '0.610'
$r = $ws.Range("E7")
$r.NumberFormat = "@"
$r.Value = '  -2.65%  '
$r = $ws.Range("E8")
$r.NumberFormat = "@"
$r.Value = '  -0.02%  '
$r = $ws.Range("E9")
$r.NumberFormat = "@"
$r.Value = '  -1.65%  '
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = '0.630'
$r = $ws.Range("E10")
$r.NumberFormat = "@"
$r.Value = '  -2.79%  '
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = '51.86'
$r = $ws.Range("E11")
$r.NumberFormat = "@"
$r.Value = '  -3.70%  '
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = '0.0000286'
$r = $ws.Range("E12")
$r.NumberFormat = "@"
$r.Value = '  -5.99%  '
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = '9.24'
$r = $ws.Range("E13")
$r.NumberFormat = "@"
$r.Value = '  -3.38%  '
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = '4.101.58'
$r = $ws.Range("E14")
$r.NumberFormat = "@"
$r.Value = '  -1.69%  '
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = '664.90'
$r = $ws.Range("E15")
$r.NumberFormat = "@"
$r.Value = '  +12.50%  '
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = '69.562.44'
$r = $ws.Range("E16")
$r.NumberFormat = "@"
$r.Value = '  -1.00%  '
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = '3.549.19'
$r = $ws.Range("E17")
$r.NumberFormat = "@"
$r.Value = '  -1.46%  '
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = '12.49'
$r = $ws.Range("E18")
$r.NumberFormat = "@"
$r.Value = '  -4.88%  '
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = '0.121'
$r = $ws.Range("E19")
$r.NumberFormat = "@"
$r.Value = '  -0.82%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '18.48'
$r = $ws.Range("E20")
$r.NumberFormat = "@"
$r.Value = '  -3.45%  '
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = '0.965'
$r = $ws.Range("E21")
$r.NumberFormat = "@"
$r.Value = '  -2.94%  '
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = '18.43'
$r = $ws.Range("E22")
$r.NumberFormat = "@"
$r.Value = '  +4.15%  '
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '5.28'
$r = $ws.Range("E23")
$r.NumberFormat = "@"
$r.Value = '  +2.33%  '
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = '105.10'
$r = $ws.Range("E24")
$r.NumberFormat = "@"
$r.Value = '  +3.34%  '
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = '4.37'
$r = $ws.Range("E25")
$r.NumberFormat = "@"
$r.Value = '  -5.05%  '
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = '2.91'
$r = $ws.Range("E26")
$r.NumberFormat = "@"
$r.Value = '  -3.85%  '
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = '10.16'
$r = $ws.Range("E27")
$r.NumberFormat = "@"
$r.Value = '  -5.58%  '
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = '9.60'
$r = $ws.Range("E28")
$r.NumberFormat = "@"
$r.Value = '  +0.06%  '
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = '33.25'
$r = $ws.Range("E29")
$r.NumberFormat = "@"
$r.Value = '  -1.95%  '
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = '4.42'
$r = $ws.Range("E30")
$r.NumberFormat = "@"
$r.Value = '  -8.02%  '
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = '6.78'
$r = $ws.Range("E31")
$r.NumberFormat = "@"
$r.Value = '  -4.70%  '
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = '11.77'
$r = $ws.Range("E32")
$r.NumberFormat = "@"
$r.Value = '  -4.28%  '
$r = $ws.Range("E33")
$r.NumberFormat = "@"
$r.Value = '  -5.17%  '
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = '61.91'
$r = $ws.Range("E34")
$r.NumberFormat = "@"
$r.Value = '  -2.19%  '
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = '3.785.46'
$r = $ws.Range("E35")
$r.NumberFormat = "@"
$r.Value = '  -4.21%  '
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = '3.75'
$r = $ws.Range("E36")
$r.NumberFormat = "@"
$r.Value = '  +6.12%  '
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = '0.0₃0815'
$r = $ws.Range("E37")
$r.NumberFormat = "@"
$r.Value = '  -8.74%  '
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '0.998'
$r = $ws.Range("E38")
$r.NumberFormat = "@"
$r.Value = '  -0.16%  '
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = '504.04'
$r = $ws.Range("E39")
$r.NumberFormat = "@"
$r.Value = '  -3.79%  '
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = '2.92'
$r = $ws.Range("E40")
$r.NumberFormat = "@"
$r.Value = '  -6.51%  '
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = '0.372'
$r = $ws.Range("E41")
$r.NumberFormat = "@"
$r.Value = '  -5.05%  '
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = '0.134'
$r = $ws.Range("E42")
$r.NumberFormat = "@"
$r.Value = '  +0.17%  '
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '34.49'
$r = $ws.Range("E43")
$r.NumberFormat = "@"
$r.Value = '  -6.51%  '
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = '0.0448'
$r = $ws.Range("E44")
$r.NumberFormat = "@"
$r.Value = '  -1.52%  '
$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '3.39'
$r = $ws.Range("E45")
$r.NumberFormat = "@"
$r.Value = '  -1.06%  '
$ws.Range("B46").Value = 'ThetaToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = '2.86'
$r = $ws.Range("E46")
$r.NumberFormat = "@"
$r.Value = '  +0.44%  '
$r = $ws.Range("E47")
$r.NumberFormat = "@"
$r.Value = '  -2.94%  '
$r = $ws.Range("E48")
$r.NumberFormat = "@"
$r.Value = '  -0.21%  '
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = '8.29'
$r = $ws.Range("E49")
$r.NumberFormat = "@"
$r.Value = '  -3.78%  '
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = '1.80'
$r = $ws.Range("E50")
$r.NumberFormat = "@"
$r.Value = '  +21.35%  '
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = '2.67'
$r = $ws.Range("E51")
$r.NumberFormat = "@"
$r.Value = '  +63.47%  '
